$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.350.68"
$ws.Range("E2").Value = "  -2.67%  "

# Row 3
$ws.Range("D3").Value = "2.220.74"
$ws.Range("E3").Value = "  -1.87%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "108.10"
$ws.Range("E5").Value = "  -8.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.89"
$ws.Range("E6").Value = "  +12.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -2.99%  "

# Row 8
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  -3.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.72"
$ws.Range("E10").Value = "  -7.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -3.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.47"
$ws.Range("E12").Value = "  +0.44%  "

# Row 13
$ws.Range("E13").Value = "  -4.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.992"
$ws.Range("E14").Value = "  +9.84%  "

# Row 15
$ws.Range("E15").Value = "  -2.45%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.98"
$ws.Range("E16").Value = "  -1.97%  "

# Row 17
$ws.Range("D17").Value = "2.550.82"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18
$ws.Range("D18").Value = "2.233.50"
$ws.Range("E18").Value = "  -1.26%  "

# Row 19
$ws.Range("D19").Value = "42.282.96"
$ws.Range("E19").Value = "  -2.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  +7.88%  "

# Row 21
$ws.Range("E21").Value = "  -3.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.44"
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +22.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -3.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "228.48"
$ws.Range("E25").Value = "  -2.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.11"
$ws.Range("E26").Value = "  -3.99%  "

# Row 27
$ws.Range("E27").Value = "  -1.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.58"
$ws.Range("E28").Value = "  -2.16%  "

# Row 29
$ws.Range("E29").Value = "  -0.84%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.17"
$ws.Range("E30").Value = "  -7.60%  "

# Row 31
$ws.Range("E31").Value = "  -4.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.94"
$ws.Range("E32").Value = "  +1.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.96"
$ws.Range("E33").Value = "  -3.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0898"
$ws.Range("E34").Value = "  -1.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.60"
$ws.Range("E35").Value = "  -1.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.08"
$ws.Range("E36").Value = "  +11.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.37"
$ws.Range("E37").Value = "  +1.75%  "

# Row 38
$ws.Range("E38").Value = "  -2.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0373"
$ws.Range("E39").Value = "  -0.63%  "

# Row 40
$ws.Range("E40").Value = "  -3.90%  "

# Row 41
$ws.Range("E41").Value = "  -4.61%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.96"
$ws.Range("E42").Value = "  -2.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.232"
$ws.Range("E43").Value = "  -1.69%  "

# Row 44
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.62"
$ws.Range("E45").Value = "  -8.99%  "

# Row 46
$ws.Range("E46").Value = "  -4.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.41"
$ws.Range("E47").Value = "  -6.16%  "

# Row 48
$ws.Range("E48").Value = "  +5.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.35"
$ws.Range("E49").Value = "  +2.36%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  +7.72%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.40"
$ws.Range("E51").Value = "  -1.14%  "
